# Apply the commit's edits to ValueSet-VaccineCodesCvxMvx.xlsx
$wb = $excel.ActiveWorkbook

# --- Rename the two "Include from ..." sheets -------------------------------
$wb.Worksheets.Item(2).Name = "Include #0"
$wb.Worksheets.Item(3).Name = "Include #1"

# --- Update the Metadata sheet ----------------------------------------------
$ws = $wb.Worksheets.Item(1)

# URL: pythia -> cicada
$ws.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/VaccineCodesCvxMvx"

# Date: regenerated timestamp
$ws.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws.Rows.Item(11).Insert()

# Restore the bordered "data row" formatting on the newly inserted row by
# copying the format from the row above it (Contact row).
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# The old "Immutable" / "BooleanType[null]" row (now row 15) is replaced by
# a new "Codes" / "All codes" row.
$ws.Range("A15").Value = "Codes"
$ws.Range("B15").Value = "All codes"
